# Update the weekly Fruta/Hortaliza data: reshuffle existing rows and
# append one new row (date 2023-06-07 / serial 45084) of Ciboulette data
# for "Terminal La Palmera de La Serena".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: @(row, FechaSerial, Volumen(J), PrecioMinimo(K), PrecioMaximo(L), PrecioPromedioPonderado(M), PrecioKg(P))
$data = @(
    @(2,44978,1000,1800,2000,1900,633),
    @(3,45006,1100,2000,2500,2250,750),
    @(4,45007,1160,2000,2500,2250,750),
    @(5,44911,700,1800,2000,1900,633),
    @(6,45035,1100,2000,2500,2250,750),
    @(7,44964,1000,2000,2500,2250,750),
    @(8,44965,1120,2000,2500,2250,750),
    @(9,45070,800,2000,2500,2250,750),
    @(10,45013,1100,2000,2500,2250,750),
    @(11,44910,1000,1800,2000,1900,633),
    @(12,45077,760,2000,2500,2250,750),
    @(13,45062,1100,2000,2500,2250,750),
    @(14,44999,1100,2000,2500,2250,750),
    @(15,44881,500,1900,2000,1950,650),
    @(16,44985,1000,2000,2500,2250,750),
    @(17,44883,500,1800,2000,1900,633),
    @(18,44992,1040,2000,2500,2250,750),
    @(19,45020,1200,2000,2500,2250,750),
    @(20,44970,800,2000,2500,2250,750),
    @(21,44951,800,2000,2500,2250,750),
    @(22,45084,900,2000,2500,2250,750),
    @(23,44953,1000,2000,2500,2250,750),
    @(24,44685,400,1500,2000,1750,583),
    @(25,45028,1000,2000,2500,2250,750),
    @(26,44848,1000,1500,2000,1750,583),
    @(27,44827,1200,2000,2500,2250,750),
    @(28,44971,1000,2000,2500,2250,750),
    @(29,45041,1160,2000,2500,2250,750),
    @(30,45034,1100,2000,2500,2250,750)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 4).Value = $row[1]
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 10).Value = $row[2]
    $ws.Cells.Item($r, 11).Value = $row[3]
    $ws.Cells.Item($r, 12).Value = $row[4]
    $ws.Cells.Item($r, 13).Value = $row[5]
    $ws.Cells.Item($r, 16).Value = $row[6]
}

# Row 30 is brand new - fill in the columns that are constant across every
# other row in this sheet (same market/category/unit/origin/etc.)
$ws.Cells.Item(30, 1).Value = 8
$ws.Cells.Item(30, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(30, 3).Value = "Coquimbo"
$ws.Cells.Item(30, 5).Value = 4
$ws.Cells.Item(30, 6).Value = 100112039
$ws.Cells.Item(30, 7).Value = "Ciboulette"
$ws.Cells.Item(30, 8).Value = "Sin especificar"
$ws.Cells.Item(30, 9).Value = "Primera"
$ws.Cells.Item(30, 14).Value = "`$/docena de atados"
$ws.Cells.Item(30, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(30, 17).Value = 3
$ws.Cells.Item(30, 18).Value = "Hortaliza"
